$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "max" column (column C); remaining columns shift left.
$ws.Columns("C").Delete()

# Update the numeric values in column B (now holding the former "max"
# column's position) with the new prediction values.
$ws.Range("B2").Value = 26.71163595541235
$ws.Range("B3").Value = 29.13910977627147
